# Apply the metadata updates to the "Metadata" worksheet of the CodeSystem workbook:
#   - B7  (Experimental value): set to the literal text "false"
#   - B8  (Date value):         bump the generation timestamp
#   - B17 (Description value):  fill in the CodeSystem description
#
# B7 needs special handling: Excel's Range.Value setter auto-coerces a bare
# "false"/"true" string literal into a native Boolean, which would serialize
# as a <c t="b"> cell instead of the plain shared-string text cell the sheet
# already uses everywhere else. Writing it as a formula that evaluates to the
# text string, then converting that formula to its static value in place via
# Copy/PasteSpecial (values only), yields a genuine text cell without
# disturbing the cell's existing style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B7").Formula = "=""false"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"

$ws.Range("B17").Value = "Cardiorespiratory fitness classification categories"
